# Updated Azure Storage Account Map
# - Remove the "wallpapers" row entirely (old row 7, which held the
#   "(none)" Naming Convention and the long wallpaper description).
# - The "activitylogs" row (old row 8) shifts up to row 7, and its
#   Naming Convention text is corrected from dashes to dots:
#   YYYY.MM.DD (i.e. "2020-8-2", 2020-8-16") -> YYYY.MM.DD (i.e. "2020.8.2", 2020.8.16")
# - The "messagesubmissions" row (old row 9) shifts up to row 8, and its
#   Naming Convention text changes from "YYYY.MM.DD.HH.MM.SS" to "YYYY.MM.DD".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Containers")

# Old layout (rows 1..9):
#   7: Admin | wallpapers          | <wallpaper description>                | (none)
#   8: Admin | activitylogs        | <activity log description>             | YYYY.MM.DD (i.e. "2020-8-2", 2020-8-16")
#   9: Admin | messagesubmissions  | Contains a series of messages to owner. | YYYY.MM.DD.HH.MM.SS
#
# New layout (rows 1..8):
#   7: Admin | activitylogs        | <activity log description>             | YYYY.MM.DD (i.e. "2020.8.2", 2020.8.16")
#   8: Admin | messagesubmissions  | Contains a series of messages to owner. | YYYY.MM.DD

# Delete the wallpapers row (row 7); this shifts rows 8-9 (activitylogs,
# messagesubmissions) up to 7-8, preserving their existing Container Name
# and Description text automatically.
$ws.Rows.Item(7).Delete()

# Fix up the Naming Convention text for the (now) row 8 - messagesubmissions
# first, so new shared-string entries are appended in the same order as the
# target workbook (messagesubmissions' short convention before the longer
# activitylogs one).
$ws.Range("D8").Value = "YYYY.MM.DD"

# Fix up the Naming Convention text for the (now) row 7 - activitylogs.
$ws.Range("D7").Value = 'YYYY.MM.DD (i.e. "2020.8.2", 2020.8.16")'

# Move the active selection to A7, matching the post-edit saved state.
$ws.Range("A7").Select()
